# Updates the crypto price/volume table on Sheet1 (rows 2-51) to the
# latest scraped values. Column D = Price, Column E = Volume(1h).
#
# Several Price values are plain decimal numbers (e.g. "0.628", "244.44").
# Excel's COM type-inference would silently coerce a bare Range.Value
# assignment like that into a numeric cell, losing the original "text"
# storage used by the source data (multi-dot numbers, subscript digits,
# trailing zeros, etc. all need to stay literal text). To keep those cells
# text after the write, NumberFormat is forced to "@" immediately before
# assigning any Price value that would otherwise look like a pure number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.247.70"
$ws.Range("E2").Value = "  -5.97%  "
$ws.Range("D3").Value = "2.218.68"
$ws.Range("E3").Value = "  -6.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.44"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -6.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.59"
$ws.Range("E7").Value = "  -6.47%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -7.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.18"
$ws.Range("E10").Value = "  +5.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0953"
$ws.Range("E11").Value = "  -7.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.18"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.70"
$ws.Range("E14").Value = "  -7.95%  "
$ws.Range("D15").Value = "2.550.29"
$ws.Range("E15").Value = "  -6.22%  "
$ws.Range("E16").Value = "  -10.16%  "
$ws.Range("E17").Value = "  -9.32%  "
$ws.Range("D18").Value = "2.221.92"
$ws.Range("E18").Value = "  -5.62%  "
$ws.Range("D19").Value = "41.294.09"
$ws.Range("E19").Value = "  -5.78%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -8.27%  "
$ws.Range("E21").Value = "  -6.52%  "
$ws.Range("E22").Value = "  -7.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.88"
$ws.Range("E23").Value = "  -8.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  +12.89%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -4.63%  "
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("E28").Value = "  -7.77%  "
$ws.Range("E29").Value = "  -5.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.82"
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("E31").Value = "  -8.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  -7.96%  "
$ws.Range("E33").Value = "  -7.08%  "
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("E36").Value = "  -9.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.90"
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.11"
$ws.Range("E38").Value = "  +16.31%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  -4.84%  "
$ws.Range("E41").Value = "  -11.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.94"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("E43").Value = "  -11.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.206"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.85"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("E47").Value = "  +11.24%  "
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("E50").Value = "  -5.85%  "
$ws.Range("E51").Value = "  -5.58%  "
